$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4; existing rows 4..29 shift down to 5..30
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly record
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C4").Value = "Metropolitana"
$ws.Range("D4").Value = 44847
$ws.Range("E4").Value = 13
$ws.Range("F4").Value = 300000001
$ws.Range("G4").Value = "Rabanito"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 7900
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("N4").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O4").Value = "Provincia de Chacabuco"
$ws.Range("P4").Value = 30
$ws.Range("Q4").Value = 100
$ws.Range("R4").Value = "Hortaliza"

# Make sure the date cell keeps the same date number format as the rest of column D
$ws.Range("D4").NumberFormat = $ws.Range("D5").NumberFormat

Write-Host "Applied weekly insert; new dimension: $($ws.UsedRange.Address())"
